$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.256.74"
$ws.Range("E2").Value = "  -1.68%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.822.26"
$ws.Range("E3").Value = "  -1.99%  "

$ws.Range("E4").Value = "  -1.29%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.41"
$ws.Range("E5").Value = "  -2.03%  "

$ws.Range("E6").Value = "  -1.28%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4270"
$ws.Range("E7").Value = "  -2.30%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3672"
$ws.Range("E8").Value = "  -2.79%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.91"
$ws.Range("E9").Value = "  -2.29%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07231"
$ws.Range("E10").Value = "  -2.52%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8609"
$ws.Range("E11").Value = "  -2.61%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.96"
$ws.Range("E12").Value = "  -2.84%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.842.27"
$ws.Range("E13").Value = "  -1.14%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.665"
$ws.Range("E14").Value = "  -1.38%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07103"
$ws.Range("E15").Value = "  -0.52%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.303"
$ws.Range("E16").Value = "  -3.45%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.28"
$ws.Range("E17").Value = "  +0.45%  "

$ws.Range("E18").Value = "  -1.60%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008865"
$ws.Range("E19").Value = "  -1.90%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.005"
$ws.Range("E20").Value = "  -1.27%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.03"
$ws.Range("E21").Value = "  -2.49%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "27.277.39"
$ws.Range("E22").Value = "  -1.60%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.138"
$ws.Range("E23").Value = "  -2.78%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.89"
$ws.Range("E24").Value = "  -2.41%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.054.98"
$ws.Range("E25").Value = "  -2.32%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.003"
$ws.Range("E26").Value = "  -1.59%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "153.25"
$ws.Range("E27").Value = "  -2.51%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.30"
$ws.Range("E28").Value = "  -2.08%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.116"
$ws.Range("E29").Value = "  +6.24%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.221"
$ws.Range("E30").Value = "  -3.86%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "116.17"
$ws.Range("E31").Value = "  -4.12%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08889"
$ws.Range("E32").Value = "  -1.80%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.195"
$ws.Range("E33").Value = "  -1.94%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7583"
$ws.Range("E34").Value = "  -1.61%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.444"
$ws.Range("E35").Value = "  -2.62%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.823"
$ws.Range("E36").Value = "  -7.04%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.005"
$ws.Range("E37").Value = "  -1.27%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.114"
$ws.Range("E38").Value = "  -2.19%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01964"
$ws.Range("E39").Value = "  -0.91%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05266"
$ws.Range("E40").Value = "  -0.79%  "

$ws.Range("E41").Value = "  +1.10%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.125"
$ws.Range("E42").Value = "  +2.31%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1677"
$ws.Range("E43").Value = "  -0.05%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5029"
$ws.Range("E44").Value = "  -2.93%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.605"
$ws.Range("E45").Value = "  -1.21%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.53"
$ws.Range("E46").Value = "  -2.42%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "106.77"
$ws.Range("E47").Value = "  -3.18%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4703"
$ws.Range("E48").Value = "  -0.64%  "

$ws.Range("E49").Value = "  -1.44%  "

$ws.Range("E50").Value = "  -1.45%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.662"
$ws.Range("E51").Value = "  -3.10%  "
